# Apply the "Updated cryptos list" data refresh to Sheet1.
# Columns: A=#, B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are text that merely look numeric (e.g. "1.439.43",
# "0.00001017"); writing them straight into .Value would make Excel coerce
# them into real numbers, so each one is entered with a leading apostrophe to
# force text, and the cell style is then reset to "Normal" so no stray
# number-format/quote-prefix style sticks around on the cell.
function Set-PriceText {
    param($range, [string]$text)
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# Column E (Volume(1h)) values already contain spaces/percent signs, so Excel
# keeps them as plain text and they can be assigned directly.
function Set-PlainText {
    param($range, [string]$text)
    $ws.Range($range).Value = $text
}

$priceUpdates = @{
    "D2"  = "20.248.05"
    "D3"  = "1.439.01"
    "D4"  = "1.008"
    "D5"  = "0.9066"
    "D6"  = "277.66"
    "D7"  = "0.3652"
    "D8"  = "0.3111"
    "D10" = "1.022"
    "D11" = "0.06537"
    "D13" = "5.381"
    "D14" = "17.58"
    "D15" = "6.059"
    "D16" = "0.00001018"
    "D17" = "1.440.40"
    "D18" = "0.9420"
    "D19" = "0.05650"
    "D20" = "67.72"
    "D21" = "5.389"
    "D22" = "14.39"
    "D23" = "10.79"
    "D24" = "2.233"
    "D25" = "20.288.12"
    "D26" = "2.162"
    "D27" = "137.47"
    "D28" = "16.93"
    "D29" = "1.591.47"
    "D30" = "110.07"
    "D31" = "3.904"
    "D32" = "0.8035"
    "D33" = "4.810"
    "D34" = "0.07679"
    "D35" = "0.05932"
    "D36" = "1.437"
    "D37" = "1.141"
    "D38" = "4.642"
    "D39" = "0.01987"
    "D40" = "10.18"
    "D41" = "0.1838"
    "D42" = "0.9126"
    "D43" = "7.066"
    "D44" = "0.5229"
    "D45" = "3.513"
    "D47" = "118.72"
    "D48" = "0.5138"
    "D49" = "1.758"
    "D50" = "0.06325"
    "D51" = "0.9874"
}

$volumeUpdates = @{
    "E2"  = "  +2.46%  "
    "E3"  = "  +3.76%  "
    "E4"  = "  +0.55%  "
    "E5"  = "  -9.66%  "
    "E6"  = "  +3.43%  "
    "E7"  = "  +0.47%  "
    "E9"  = "  +2.34%  "
    "E10" = "  +4.95%  "
    "E11" = "  +2.06%  "
    "E12" = "  -0.43%  "
    "E13" = "  +1.96%  "
    "E14" = "  +7.13%  "
    "E15" = "  +0.24%  "
    "E16" = "  +2.92%  "
    "E17" = "  +3.65%  "
    "E18" = "  -6.14%  "
    "E19" = "  +0.46%  "
    "E20" = "  -2.97%  "
    "E21" = "  -2.01%  "
    "E22" = "  +0.65%  "
    "E23" = "  +2.66%  "
    "E24" = "  -0.36%  "
    "E25" = "  +2.74%  "
    "E26" = "  +0.31%  "
    "E27" = "  +1.15%  "
    "E28" = "  +2.48%  "
    "E29" = "  +2.99%  "
    "E30" = "  +2.47%  "
    "E31" = "  +2.28%  "
    "E32" = "  +2.09%  "
    "E33" = "  -7.62%  "
    "E34" = "  +1.44%  "
    "E35" = "  +6.44%  "
    "E36" = "  +11.67%  "
    "E37" = "  +9.23%  "
    "E38" = "  -0.66%  "
    "E39" = "  -1.06%  "
    "E40" = "  +2.47%  "
    "E41" = "  -1.51%  "
    "E42" = "  -8.99%  "
    "E43" = "  -13.96%  "
    "E44" = "  +1.50%  "
    "E45" = "  +1.50%  "
    "E46" = "  +2.07%  "
    "E47" = "  +10.04%  "
    "E48" = "  +3.74%  "
    "E49" = "  +2.51%  "
    "E50" = "  +4.87%  "
    "E51" = "  -1.76%  "
}

foreach ($cellRef in $priceUpdates.Keys) {
    Set-PriceText $cellRef $priceUpdates[$cellRef]
}

foreach ($cellRef in $volumeUpdates.Keys) {
    Set-PlainText $cellRef $volumeUpdates[$cellRef]
}

# Rows 44/45 also swapped places: PancakeSwap <-> TheSandbox (name + link),
# in addition to the price/volume updates already handled above.
Set-PlainText "B44" "TheSandbox"
Set-PlainText "C44" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-PlainText "B45" "PancakeSwap"
Set-PlainText "C45" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
